# Weekend-Weekday model validation over different months.
# Updates evap / Inflow data series and downstream Scalar summary values.

$wb = $excel.ActiveWorkbook

# --- evap sheet: B4:B34 values replaced with Inflow's old values ---
$wsEvap = $wb.Worksheets.Item("evap")
$evapValues = @(1732,1730,1727,1726,1723,1720,1719,1717,1714,1711,1709,1706,1703,1702,1700,1697,1695,1692,1690,1687,1686,1684,1681,1679,1676,1673,1670,1668,1666,1664,1661)
for ($i = 0; $i -lt $evapValues.Length; $i++) {
    $row = 4 + $i
    $wsEvap.Range("B$row").Value = $evapValues[$i]
}

# --- Inflow sheet: B4:B34 values replaced with new data ---
$wsInflow = $wb.Worksheets.Item("Inflow")
$inflowValues = @(2724,3660,5360,3831,4188,4340,3944,3349,3215,4257,4265,4177,3824,5149,2978,3836,5588,5659,4304,5190,6253,4172,4503,4295,2236,4424,3682,2998,3025,4101,3469)
for ($i = 0; $i -lt $inflowValues.Length; $i++) {
    $row = 4 + $i
    $wsInflow.Range("B$row").Value = $inflowValues[$i]
}

# --- Scalar sheet: downstream summary values recalculated ---
$wsScalar = $wb.Worksheets.Item("Scalar")
$wsScalar.Range("B12").Value = 12030761.436960001
$wsScalar.Range("C22").Value = 12964160.032
$wsScalar.Range("C25").Value = 12030761.436960001
$wsScalar.Range("B32").Value = -52608
$wsScalar.Range("C32").Value = -52608
$wsScalar.Range("D32").Value = -52608
$wsScalar.Range("C33").Value = 12030761.436960001
$wsScalar.Range("C34").Value = 12030761.436960001
$wsScalar.Range("B36").Value = 12964160.032
$wsScalar.Range("C36").Value = 12964160.032
$wsScalar.Range("D36").Value = 12964160.032
